$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Before touching rows, copy the "last row" bottom-border formatting
#        from row 30 (old last data row) onto row 22 (new last data row),
#        so once rows 23-30 are removed, row 22 keeps the emphasized bottom
#        border that a final table row has.
$ws.Range("B30:J30").Copy() | Out-Null
$ws.Range("B22:J22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Update the summary header values
$ws.Range("E11").Value = 506820
$ws.Range("C13").Value = 6

# --- 3. Update the data table (rows 16-22)
# Row 16: CEDRICK CONTRERA GUARDO, period changes 2504 -> 2507 (amounts unchanged)
$ws.Range("E16").Value = "2507"

# Row 17: now XIOMARA PATRICIA MONROY TINOCO, period 2507, 128000 / 3200000
$ws.Range("C17").Value = "45687489"
$ws.Range("D17").Value = "XIOMARA PATRICIA MONROY TINOCO"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 128000
$ws.Range("G17").Value = 3200000

# Row 18: now KAREN MARGARITA SARABIA AYOLA, period 2507, 56940 / 1423500
$ws.Range("C18").Value = "32939066"
$ws.Range("D18").Value = "KAREN MARGARITA SARABIA AYOLA"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19: NELSON ENRIQUE GOMEZ VEGA, period 2507 (amounts unchanged)
$ws.Range("C19").Value = "1047475016"
$ws.Range("D19").Value = "NELSON ENRIQUE GOMEZ VEGA"
$ws.Range("E19").Value = "2507"

# Row 20: DARWIN ALEXANDER ESCOBAR MIRANDA, period 2507 (amounts unchanged)
$ws.Range("C20").Value = "1047434781"
$ws.Range("D20").Value = "DARWIN ALEXANDER ESCOBAR MIRANDA"
$ws.Range("E20").Value = "2507"

# Row 21: new worker GLENDA MARCELA MIELES GOMEZ, period 2301, 40000 / 1200000
$ws.Range("C21").Value = "1002191142"
$ws.Range("D21").Value = "GLENDA MARCELA MIELES GOMEZ"
$ws.Range("E21").Value = "2301"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1200000

# Row 22: same worker GLENDA MARCELA MIELES GOMEZ, period 2212, 40000 / 1200000
$ws.Range("C22").Value = "1002191142"
$ws.Range("D22").Value = "GLENDA MARCELA MIELES GOMEZ"
$ws.Range("E22").Value = "2212"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 1200000

# --- 4. Remove the old rows 23-30 (period 2505 / 2506 data, now obsolete)
$ws.Rows("23:30").Delete()
